$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("!!Main root")
$ws1.Range("A1").Value = "!!!ObjTables objTablesVersion='0.0.8'"
$ws1.Range("A2").Value = "!!ObjTables type='Data' id='MainRoot'"

$ws2 = $wb.Worksheets.Item("!!Nodes")
$ws2.Range("A1").Value = "!!ObjTables type='Data' id='Node'"

$ws3 = $wb.Worksheets.Item("!!Leaves")
$ws3.Range("A1").Value = "!!ObjTables type='Data' id='Leaf'"

$ws4 = $wb.Worksheets.Item("!!One to many rows")
$ws4.Range("A1").Value = "!!ObjTables type='Data' id='OneToManyRow'"
